$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.543.70'
$ws.Range("E2").Value = '  -2.11%  '

$ws.Range("D3").Value = '2.409.10'
$ws.Range("E3").Value = '  -1.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.993'
$ws.Range("E4").Value = '  -0.53%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.17'
$ws.Range("E5").Value = '  -2.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.41'
$ws.Range("E6").Value = '  -2.92%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +0.85%  '

$ws.Range("D9").Value = '2.389.58'
$ws.Range("E9").Value = '  -2.16%  '

$ws.Range("E10").Value = '  -3.88%  '

$ws.Range("E11").Value = '  -1.19%  '

$ws.Range("E12").Value = '  -2.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.61'
$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("D15").Value = '2.823.07'
$ws.Range("E15").Value = '  -2.42%  '

$ws.Range("E16").Value = '  -2.75%  '

$ws.Range("D17").Value = '60.715.18'
$ws.Range("E17").Value = '  -1.68%  '

$ws.Range("D18").Value = '2.383.11'
$ws.Range("E18").Value = '  -2.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.04'
$ws.Range("E19").Value = '  +11.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.52'
$ws.Range("E20").Value = '  -0.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.18'
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("E22").Value = '  -1.03%  '

$ws.Range("E23").Value = '  -2.62%  '

$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.80'
$ws.Range("E25").Value = '  -6.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.09'
$ws.Range("E26").Value = '  -1.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.18'
$ws.Range("E27").Value = '  -10.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '557.05'
$ws.Range("E28").Value = '  -4.47%  '

$ws.Range("D29").Value = '2.517.94'
$ws.Range("E29").Value = '  -2.00%  '

$ws.Range("D30").Value = '0.0₃0917'
$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("E32").Value = '  -5.26%  '

$ws.Range("E33").Value = '  -4.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.131'
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  +0.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.31'
$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("E38").Value = '  -1.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.52'
$ws.Range("E39").Value = '  -4.67%  '

$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.06'
$ws.Range("E41").Value = '  -1.45%  '

$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("E43").Value = '  -1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.33'
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").Value = '0.0₆0289'
$ws.Range("E45").Value = '  +3.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.10'
$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.50'
$ws.Range("E47").Value = '  -1.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.581'
$ws.Range("E48").Value = '  -2.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0497'
$ws.Range("E49").Value = '  -2.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.09'
$ws.Range("E50").Value = '  -2.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0896'
$ws.Range("E51").Value = '  +0.09%  '
